$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Tue Oct 31 00:50:03 2023"
$ws.Range("B4").Value = "HKHSI"
$ws.Range("C4").Value = "Hang Seng Indexes"
$ws.Range("D4").Value = "未开盘"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "17406.36"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "+7.63  +0.04%"
$ws.Range("G4").Value = 17406.36
$ws.Range("H4").Value = 17225.65
$ws.Range("I4").Value = 22700.85
$ws.Range("J4").Value = 0.01
$ws.Range("K4").Value = 17225.65
$ws.Range("L4").Value = 17398.73
$ws.Range("M4").Value = 14597.31
$ws.Range("N4").Value = 92472000000
